# Insert a new weekly record at row 37 (Fruta / hortaliza, semanal)
# This shifts the existing rows 37:87 down to 38:88 and fills the
# newly inserted row 37 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37 (pushes rows 37..87 down to 38..88)
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value2 = 7
$ws.Range("B37").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value2 = "Ñuble"
$ws.Range("D37").Value2 = 44771
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("E37").Value2 = 16
$ws.Range("F37").Value2 = 100112031
$ws.Range("G37").Value2 = "Poroto verde"
$ws.Range("H37").Value2 = "Sin especificar"
$ws.Range("I37").Value2 = "Primera"
$ws.Range("J37").Value2 = 40
$ws.Range("K37").Value2 = 30000
$ws.Range("L37").Value2 = 30000
$ws.Range("M37").Value2 = 30000
$ws.Range("N37").Value2 = "`$/malla 25 kilos"
$ws.Range("O37").Value2 = "Región de Arica y Parinacota"
$ws.Range("P37").Value2 = 1200
$ws.Range("Q37").Value2 = 25
$ws.Range("R37").Value2 = "Hortaliza"
